# Append the "gist" note to the end of the document, after the existing
# pictures and the two trailing blank paragraphs:
#
#   ... (last picture paragraph)
#   <blank paragraph>
#   <blank paragraph>
#   "Same code is uploaded on gist as github integration on Watson studio was giving error."
#   "Gist link:"
#   "https://gist.github.com/ashitole/308beb105f64e37f520e3cdfa3539223"

$d = $word.ActiveDocument

function New-TrailingParagraph {
    # Adds a new (empty) paragraph at the very end of the document's story.
    $r = $d.Range($d.Content.End, $d.Content.End)
    $r.InsertParagraphAfter() | Out-Null
}

function Add-TrailingText([string]$text) {
    # Inserts text just before the final paragraph mark, i.e. into the
    # document's current last paragraph.
    $r = $d.Range($d.Content.End - 1, $d.Content.End - 1)
    $r.InsertAfter($text) | Out-Null
}

# Paragraph 1: "Same code is uploaded on gist as github integration on
# Watson studio was giving error." (typed as three chunks so that "github"
# lines up with its own run, matching how the word got flagged/split by
# the spell checker in the original edit).
New-TrailingParagraph
Add-TrailingText "Same code is uploaded on gist as "
Add-TrailingText "github"
Add-TrailingText " integration on Watson studio was giving error."

# Paragraph 2: "Gist link:"
New-TrailingParagraph
Add-TrailingText "Gist link:"

# Paragraph 3: the gist URL itself.
New-TrailingParagraph
Add-TrailingText "https://gist.github.com/ashitole/308beb105f64e37f520e3cdfa3539223"

Write-Host "Paragraphs after edit:" $d.Paragraphs.Count
